# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This applies the new "K" column (column G) values for rows 2-72 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for G2:G72 (column "K"), in row order.
$newKValues = @(2, 0, 2, 3, 0, 3, 1, 1, 3, 0, 0, 2, 3, 0, 0, 1, 2, 0, 0, 3, 2, 1, 2, 2, 1, 2, 0, 1, 2, 2, 0, 1, 1, 0, 1, 2, 3, 2, 1, 2, 1, 1, 2, 2, 2, 2, 1, 1, 0, 1, 1, 1, 0, 1, 2, 2, 2, 2, 0, 1, 1, 2, 2, 3, 2, 1, 1, 1, 1, 0, 1)

$startRow = 2
for ($i = 0; $i -lt $newKValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newKValues[$i]
}

$wb.Save()
